$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 14
$ws.Range("B4").Value = "<hat>"
$ws.Range("C5").Value = 15
$ws.Range("C6").Value = 16
$ws.Range("C7").Value = 11
$ws.Range("C8").Value = 16
$ws.Range("C9").Value = 15
$ws.Range("C10").Value = 12
$ws.Range("C12").Value = 17
$ws.Range("C13").Value = 13
$ws.Range("C15").Value = 16
$ws.Range("C16").Value = 15
$ws.Range("B18").Value = "<out>"
$ws.Range("C18").Value = 10
